# Insert a new weekly price record at row 19 ("Comercializadora del Agro de
# Limarí" / Zapallo italiano). This pushes all existing rows 19:48 down to
# 20:49, preserving their data, and leaves a blank row 19 to be filled in
# with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44580
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112032
$ws.Range("G19").Value = "Zapallo italiano"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12500
$ws.Range("N19").Value = "$/caja 60 unidades"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 208
$ws.Range("Q19").Value = 60
$ws.Range("R19").Value = "Hortaliza"
